$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1923.4706
$ws.Range("J17").Value = 2549.4
$ws.Range("L17").Value = 7648.200000000001
$ws.Range("N17").Value = -7984.200000000001
$ws.Range("H109").Value = 83591.664
$ws.Range("J109").Value = 83591.664
$ws.Range("L109").Value = 83591.664
$ws.Range("N109").Value = -86365.664
$ws.Range("H113").Value = 10334.429
$ws.Range("I113").Value = 23318.6
$ws.Range("J113").Value = 3121
$ws.Range("K113").Value = 23318.6
$ws.Range("L113").Value = 3121
$ws.Range("M113").Value = -20064.6
$ws.Range("N113").Value = -9629
$ws.Range("H116").Value = 1196413.5
$ws.Range("I116").Value = 5498.5
$ws.Range("K116").Value = 5498.5
$ws.Range("M116").Value = -2056.5
$ws.Range("H132").Value = 1813.9025
$ws.Range("I132").Value = 1276.5883
$ws.Range("K132").Value = 3829.7649
$ws.Range("M132").Value = -1299.7649
$ws.Range("H133").Value = 86963.75
$ws.Range("J133").Value = 86963.75
$ws.Range("L133").Value = 86963.75
$ws.Range("N133").Value = -97083.75
$ws.Range("H134").Value = 52135.715
$ws.Range("J134").Value = 55491.668
$ws.Range("L134").Value = 55491.668
$ws.Range("N134").Value = -65631.66800000001
$ws.Range("H136").Value = 96491.664
$ws.Range("J136").Value = 96491.664
$ws.Range("L136").Value = 96491.664
$ws.Range("N136").Value = -106691.664
$ws.Range("H137").Value = 764739.4
$ws.Range("I137").Value = 1415.8182
$ws.Range("J137").Value = 1814309.2
$ws.Range("K137").Value = 4247.4546
$ws.Range("L137").Value = 5442927.6
$ws.Range("M137").Value = -1697.4546
$ws.Range("N137").Value = -5448027.6
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
$ws.Range("H140").Value = 91491.664
$ws.Range("J140").Value = 91990
$ws.Range("L140").Value = 91990
$ws.Range("N140").Value = -102350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7517.73
$ws.Range("I32").Value = 4382.9644
$ws.Range("K32").Value = 4382.9644
$ws.Range("M32").Value = -4095.9644
$ws.Range("H61").Value = 54631.316
$ws.Range("I61").Value = 1706.4667
$ws.Range("K61").Value = 1706.4667
$ws.Range("M61").Value = -1494.4667
$ws.Range("H136").Value = 54631.316
$ws.Range("I136").Value = 1706.4667
$ws.Range("K136").Value = 5119.4001
$ws.Range("M136").Value = -2569.4001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 85626.914
$ws.Range("I20").Value = 144506.28
$ws.Range("K20").Value = 144506.28
$ws.Range("M20").Value = -144259.28
$ws.Range("H99").Value = 2937684
$ws.Range("I99").Value = 132253.25
$ws.Range("J99").Value = 10418833
$ws.Range("K99").Value = 132253.25
$ws.Range("L99").Value = 10418833
$ws.Range("M99").Value = -130755.25
$ws.Range("N99").Value = -10421829
$ws.Range("H107").Value = 4354.467
$ws.Range("I107").Value = 3356.889
$ws.Range("K107").Value = 3356.889
$ws.Range("M107").Value = -1436.889
$ws.Range("H109").Value = 99988.336
$ws.Range("J109").Value = 99988.336
$ws.Range("L109").Value = 99988.336
$ws.Range("N109").Value = -102762.336
$ws.Range("H132").Value = 81272.5
$ws.Range("J132").Value = 81272.5
$ws.Range("L132").Value = 81272.5
$ws.Range("N132").Value = -91392.5
$ws.Range("H138").Value = 96368
$ws.Range("J138").Value = 96368
$ws.Range("L138").Value = 96368
$ws.Range("N138").Value = -106648
$ws.Range("H140").Value = 95118.336
$ws.Range("J140").Value = 95118.336
$ws.Range("L140").Value = 95118.336
$ws.Range("N140").Value = -105478.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1624.6666
$ws.Range("I16").Value = 1342.4736
$ws.Range("J16").Value = 2697
$ws.Range("K16").Value = 1342.4736
$ws.Range("L16").Value = 2697
$ws.Range("M16").Value = -1055.4736
$ws.Range("N16").Value = -3271
$ws.Range("H31").Value = 3460.6667
$ws.Range("I31").Value = 2426.5386
$ws.Range("J31").Value = 6149.4
$ws.Range("K31").Value = 2426.5386
$ws.Range("L31").Value = 6149.4
$ws.Range("M31").Value = -2131.5386
$ws.Range("N31").Value = -6739.4
$ws.Range("H34").Value = 3460.6667
$ws.Range("I34").Value = 2426.5386
$ws.Range("J34").Value = 6149.4
$ws.Range("K34").Value = 2426.5386
$ws.Range("L34").Value = 6149.4
$ws.Range("M34").Value = -2224.5386
$ws.Range("N34").Value = -6553.4
$ws.Range("H99").Value = 1840845.4
$ws.Range("J99").Value = 5210879
$ws.Range("L99").Value = 5210879
$ws.Range("N99").Value = -5213875
$ws.Range("H107").Value = 810.875
$ws.Range("I107").Value = 666.1667
$ws.Range("K107").Value = 666.1667
$ws.Range("M107").Value = 1253.8333
$ws.Range("H113").Value = 1624.6666
$ws.Range("I113").Value = 1342.4736
$ws.Range("J113").Value = 2697
$ws.Range("K113").Value = 1342.4736
$ws.Range("L113").Value = 2697
$ws.Range("M113").Value = 827.5264
$ws.Range("N113").Value = -7037
$ws.Range("H126").Value = 1840845.4
$ws.Range("J126").Value = 5210879
$ws.Range("L126").Value = 15632637
$ws.Range("N126").Value = -15637577
$ws.Range("H138").Value = 80050.875
$ws.Range("J138").Value = 80050.875
$ws.Range("L138").Value = 80050.875
$ws.Range("N138").Value = -90330.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224
$ws.Range("H131").Value = 34926.7
$ws.Range("I131").Value = 200675.8
$ws.Range("J131").Value = 1776.88
$ws.Range("K131").Value = 602027.3999999999
$ws.Range("L131").Value = 5330.64
$ws.Range("M131").Value = -596987.3999999999
$ws.Range("N131").Value = -15410.64
$ws.Range("H137").Value = 13710.75
$ws.Range("J137").Value = 13833.714
$ws.Range("L137").Value = 41501.142
$ws.Range("N137").Value = -51701.142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 833.625
$ws.Range("I102").Value = 630.1111
$ws.Range("K102").Value = 630.1111
$ws.Range("M102").Value = 991.8889
$ws.Range("H113").Value = 3032085.2
$ws.Range("I113").Value = 1120
$ws.Range("J113").Value = 5557889.5
$ws.Range("K113").Value = 1120
$ws.Range("L113").Value = 5557889.5
$ws.Range("M113").Value = 1050
$ws.Range("N113").Value = -5562229.5
$ws.Range("H126").Value = 3786.5386
$ws.Range("I126").Value = 2982.5
$ws.Range("J126").Value = 4475.7144
$ws.Range("K126").Value = 8947.5
$ws.Range("L126").Value = 13427.1432
$ws.Range("M126").Value = -6477.5
$ws.Range("N126").Value = -18367.1432
$ws.Range("H135").Value = 95121.664
$ws.Range("J135").Value = 95121.664
$ws.Range("L135").Value = 95121.664
$ws.Range("N135").Value = -105261.664
$ws.Range("H140").Value = 98991.664
$ws.Range("J140").Value = 98990
$ws.Range("L140").Value = 98990
$ws.Range("N140").Value = -109350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 61899.375
$ws.Range("I7").Value = 32170.715
$ws.Range("K7").Value = 32170.715
$ws.Range("M7").Value = -32058.715
$ws.Range("H40").Value = 27780280
$ws.Range("I40").Value = 5004
$ws.Range("K40").Value = 5004
$ws.Range("M40").Value = -4868
$ws.Range("H117").Value = 59193
$ws.Range("J117").Value = 59193
$ws.Range("L117").Value = 59193
$ws.Range("N117").Value = -68371
$ws.Range("H126").Value = 61899.375
$ws.Range("I126").Value = 32170.715
$ws.Range("K126").Value = 96512.145
$ws.Range("M126").Value = -94042.145
$ws.Range("H127").Value = 58120.57
$ws.Range("J127").Value = 58120.57
$ws.Range("L127").Value = 58120.57
$ws.Range("N127").Value = -68040.57000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 349.16
$ws.Range("I113").Value = 324.5
$ws.Range("J113").Value = 412.57144
$ws.Range("K113").Value = 973.5
$ws.Range("L113").Value = 1237.71432
$ws.Range("M113").Value = 1196.5
$ws.Range("N113").Value = -5577.71432
$ws.Range("H126").Value = 1751
$ws.Range("I126").Value = 1282.5625
$ws.Range("K126").Value = 3847.6875
$ws.Range("M126").Value = -1377.6875
